$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended after the existing data (rows 256-269),
# continuing the date series in column A (serial dates 44330-44343,
# i.e. 2021-05-14 through 2021-05-27) with zeros in columns B, C, D.

$startSerial = 44330
$startRow = 256
$endRow = 269

$srcDateCell = $ws.Range("A255")

for ($row = $startRow; $row -le $endRow; $row++) {
    $serial = $startSerial + ($row - $startRow)

    $dateCell = $ws.Cells.Item($row, 1)
    # Copy formatting (style) from the last existing date cell so the new
    # cell reuses the same style (border/font/alignment/number format).
    $srcDateCell.Copy($dateCell)
    $dateCell.Value2 = $serial

    $ws.Cells.Item($row, 2).Value2 = 0
    $ws.Cells.Item($row, 3).Value2 = 0
    $ws.Cells.Item($row, 4).Value2 = 0
}
